# This script updates the "cryptos" price-list worksheet to reflect the
# latest GitHub Actions scrape: most rows only get an updated Price (col D)
# value, while three coins (KickToken, BKEXToken, CEJI) swap rank-order
# positions (rows 41-43) bringing along their own Coin/Link/Price/Volume data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellAddress, $Text)
    # Force Excel to store the value as text (matching the workbook's existing
    # inline-string cells) rather than silently converting numeric-looking
    # strings into real numbers (which would also mangle formatting such as
    # trailing zeros or tiny magnitudes rendered in scientific notation).
    $Sheet.Range($CellAddress).Value = "'" + $Text
    # Clear the "quote prefix" styling side effect introduced above so the
    # cell's style index stays the same as before the edit.
    $Sheet.Range($CellAddress).Style = "Normal"
}

# ---- Simple Price-only (column D) updates -------------------------------
$priceUpdates = @{
    "D2"  = "281.69"
    "D3"  = "20.68"
    "D4"  = "6.256"
    "D5"  = "0.06145"
    "D6"  = "3.578"
    "D7"  = "6.562"
    "D8"  = "1.499"
    "D9"  = "0.8165"
    "D12" = "0.08332"
    "D13" = "0.03543"
    "D14" = "0.03181"
    "D15" = "0.09136"
    "D16" = "3.705"
    "D17" = "0.001642"
    "D18" = "0.04683"
    "D19" = "0.006535"
    "D20" = "0.006162"
    "D23" = "3.785"
    "D25" = "0.3358"
    "D40" = "0.04677"
    "D44" = "0.01112"
    "D45" = "0.00005833"
    "D48" = "0.002946"
}

foreach ($cellAddress in $priceUpdates.Keys) {
    Set-TextValue $ws $cellAddress $priceUpdates[$cellAddress]
}

# ---- Rows 41-43 rotate: CEJI, KickToken, BKEXToken shift ranking ---------
# Before: 41=KickToken, 42=BKEXToken, 43=CEJI
# After:  41=CEJI,      42=KickToken, 43=BKEXToken
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D41" "0.005604"
$ws.Range("E41").Value = "40CEJICEJI"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D42" "0.007172"
$ws.Range("E42").Value = "41KickTokenKICK"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1097"
$ws.Range("E43").Value = "42BKEXTokenBKK"
